# Generate Report for Handback
#
# - Flips the "Ready for handoff" status to "Handed back: in sync with en-US"
#   everywhere it shows up (Overview sheet's per-locale status columns, and each
#   locale sheet's own Status column).
# - Populates the previously-empty "Latest Target File" (F) / "Latest Handback
#   File" (G) columns on each locale sheet, with hyperlinks, for both data rows.
# - Stamps a real "Latest Handback DateTime" (H) instead of the zero-date
#   placeholder, per locale.
#
# Hyperlinks.Add() always appends new entries after whatever was already on the
# sheet, so to land the new F/G links in-place (between D and the next row's A)
# every hyperlink on a locale sheet is captured, cleared, and re-added in the
# desired final order. Hyperlinks.Add() also likes to stamp its own "Hyperlink"
# cell style, so the underline/blue-link formatting is re-applied to every
# linked cell afterwards to keep them visually consistent.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"
$linkColor = 15570276   # BGR-packed RGB(0x64,0x95,0xED) -> matches the workbook's HyperLink style

# ---------------------------------------------------------------------------
# Overview sheet: per-locale status columns (B = zh-cn, C = de-de)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------------
# Per-locale sheets.
# ---------------------------------------------------------------------------
$locales = @(
    @{
        Sheet        = "zh-cn"
        Tag          = "zh-cn"
        HandoffHash  = "78029f5ae1863a5a4315e72ff79149e88c6441a2"
        HandbackDate = "2016-03-20 20:27:22"
    },
    @{
        Sheet        = "de-de"
        Tag          = "de-de"
        HandoffHash  = "51321665aa757e25303718a4e43fccd2af391e1e"
        HandbackDate = "2016-03-20 20:27:28"
    }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # Latest Handback DateTime (H) - real timestamp instead of the zero-date placeholder.
    $ws.Range("H2").Value = $loc.HandbackDate
    $ws.Range("H3").Value = $loc.HandbackDate

    # Capture the 6 pre-existing hyperlinks (A2,B2,D2,A3,B3,D3) so they can be
    # re-created in the right order alongside the new F/G ones.
    $addrs = @()
    $refs = @()
    $disps = @()
    foreach ($h in $ws.Hyperlinks) {
        $addrs += $h.Address
        $refs += $h.Range.Address()
        $disps += $h.TextToDisplay
    }
    $ws.Hyperlinks.Delete()

    $xlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2." + $loc.Tag + ".xlf"
    $targetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a14241486b4eb6be7a398cdb541269aecc3a4e82/e2e/a.md"
    $handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $loc.HandoffHash + "/ol-handback/OpenLocalizationTestOrg/oltest." + $loc.Tag + "/ci/hb/" + $xlfName

    $ws.Range("F2").Value = "a.md"
    $ws.Range("G2").Value = $xlfName
    $ws.Range("F3").Value = "a.md"
    $ws.Range("G3").Value = $xlfName

    # Re-add in final order: A2,B2,D2,F2,G2,A3,B3,D3,F3,G3
    $ws.Hyperlinks.Add($ws.Range($refs[0]), $addrs[0], "", "", $disps[0])
    $ws.Hyperlinks.Add($ws.Range($refs[1]), $addrs[1], "", "", $disps[1])
    $ws.Hyperlinks.Add($ws.Range($refs[2]), $addrs[2], "", "", $disps[2])
    $ws.Hyperlinks.Add($ws.Range("F2"), $targetUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("G2"), $handbackUrl, "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range($refs[3]), $addrs[3], "", "", $disps[3])
    $ws.Hyperlinks.Add($ws.Range($refs[4]), $addrs[4], "", "", $disps[4])
    $ws.Hyperlinks.Add($ws.Range($refs[5]), $addrs[5], "", "", $disps[5])
    $ws.Hyperlinks.Add($ws.Range("F3"), $targetUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("G3"), $handbackUrl, "", "", $xlfName)

    # Keep every hyperlinked cell's look consistent (underline + the workbook's
    # custom blue), regardless of which "Hyperlink" style variant Add() stamped.
    foreach ($addr in @("A2", "B2", "D2", "F2", "G2", "A3", "B3", "D3", "F3", "G3")) {
        $ws.Range($addr).Font.Underline = $true
        $ws.Range($addr).Font.Color = $linkColor
    }
}
